$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("AREPD", 1, 6.129342789070967, 0.000000003313112006964047, "Sí"),
    @("AREPD", 2, 5.689178404966309, 0.000000040082434704658, "Sí"),
    @("AREPD", 3, 7.313096946548476, 0.000000000001692868067948439, "Sí"),
    @("AREPD", 4, 7.998272624693561, 0.00000000000001176836406102666, "Sí"),
    @("AREPD", 5, 7.570115563891723, 0.0000000000002757793993168889, "Sí"),
    @("AREPD", 6, 7.346963482349071, 0.000000000001337374655463464, "Sí"),
    @("AREPD", 7, 9.225262683519684, 0, "Sí"),
    @("AREPD", 10, 7.881894792264741, 0.00000000000002819966482547898, "Sí"),
    @("AV-MCPS", 1, 4.336048293551071, 0.0000285649688822609, "Sí"),
    @("AV-MCPS", 2, 6.022638851122547, 0.000000006161926524583805, "Sí"),
    @("AV-MCPS", 3, 6.985440623045152, 0.00000000001567812546454661, "Sí"),
    @("AV-MCPS", 4, 8.163139912905695, 0.00000000000000333066907387547, "Sí"),
    @("AV-MCPS", 5, 7.479719904530153, 0.0000000000005255795798575491, "Sí"),
    @("AV-MCPS", 6, 7.222290653849736, 0.000000000003168132423070347, "Sí"),
    @("AV-MCPS", 7, 9.29902607666186, 0, "Sí"),
    @("AV-MCPS", 10, 7.833650067132137, 0.00000000000004019007349143067, "Sí"),
    @("Block Bootstrapping", 1, 6.276460087853002, 0.000000001384601544529573, "Sí"),
    @("Block Bootstrapping", 2, 5.677138882338659, 0.00000004280528442812681, "Sí"),
    @("Block Bootstrapping", 3, 7.341522768002342, 0.000000000001389111048410996, "Sí"),
    @("Block Bootstrapping", 4, 8.0016553138109, 0.0000000000000113242748511766, "Sí"),
    @("Block Bootstrapping", 5, 7.573149321785887, 0.000000000000269784194983913, "Sí"),
    @("Block Bootstrapping", 6, 7.366605527493524, 0.000000000001165956220461339, "Sí"),
    @("Block Bootstrapping", 7, 9.220794717729014, 0, "Sí"),
    @("Block Bootstrapping", 10, 7.883528624908579, 0.00000000000002775557561562891, "Sí"),
    @("DeepAR", 1, 2.726268478861579, 0.008512010379494672, "Sí"),
    @("DeepAR", 2, 6.130772376731351, 0.000000003285452576662351, "Sí"),
    @("DeepAR", 3, 7.136724555144233, 0.000000000005679012815562601, "Sí"),
    @("DeepAR", 4, 7.501160119975993, 0.0000000000004514166818125886, "Sí"),
    @("DeepAR", 5, 7.46087808083427, 0.0000000000006008527009271347, "Sí"),
    @("DeepAR", 6, 7.154174598895085, 0.000000000005044409334686861, "Sí"),
    @("DeepAR", 7, 9.279961535812538, 0, "Sí"),
    @("DeepAR", 10, 7.938335721394182, 0.0000000000000184297022087776, "Sí"),
    @("EnCQR-LSTM", 1, 8.058656915660206, 0.000000000000007327471962526033, "Sí"),
    @("EnCQR-LSTM", 2, 5.508571626316348, 0.0000001059631213884416, "Sí"),
    @("EnCQR-LSTM", 3, 7.008497762734147, 0.00000000001344813149728452, "Sí"),
    @("EnCQR-LSTM", 4, 8.435933456084062, 0.0000000000000004440892098500626, "Sí"),
    @("EnCQR-LSTM", 5, 7.042736789615304, 0.00000000001069877519910278, "Sí"),
    @("EnCQR-LSTM", 6, 7.070293242394531, 0.000000000008893108471852429, "Sí"),
    @("EnCQR-LSTM", 7, 9.103313766460255, 0, "Sí"),
    @("EnCQR-LSTM", 10, 7.74343899946073, 0.00000000000007838174553853605, "Sí"),
    @("LSPM", 1, 11.92175987445959, 0, "Sí"),
    @("LSPM", 2, 5.792258091600652, 0.00000002271116272822837, "Sí"),
    @("LSPM", 3, 7.228751732411254, 0.000000000003030686812621752, "Sí"),
    @("LSPM", 4, 7.926553348350791, 0.00000000000002020605904817785, "Sí"),
    @("LSPM", 5, 7.361395867574009, 0.000000000001209254918421721, "Sí"),
    @("LSPM", 6, 7.169975360955958, 0.000000000004530154029680489, "Sí"),
    @("LSPM", 7, 8.994142721158154, 0, "Sí"),
    @("LSPM", 10, 8.04517649217046, 0.000000000000008215650382226158, "Sí"),
    @("LSPMW", 1, 7.927429689444206, 0.00000000000001998401444325282, "Sí"),
    @("LSPMW", 2, 5.743603937413681, 0.00000002973109003434615, "Sí"),
    @("LSPMW", 3, 7.008055413806487, 0.0000000000134878774815661, "Sí"),
    @("LSPMW", 4, 7.842391644410532, 0.00000000000003774758283725532, "Sí"),
    @("LSPMW", 5, 7.462006109410835, 0.000000000000595967719618784, "Sí"),
    @("LSPMW", 6, 7.198273973997487, 0.000000000003734568210234102, "Sí"),
    @("LSPMW", 7, 9.249628338053698, 0, "Sí"),
    @("LSPMW", 10, 7.807997874361709, 0.00000000000004862776847858186, "Sí"),
    @("MCPS", 1, 4.874556443214736, 0.000002547181449541469, "Sí"),
    @("MCPS", 2, 5.628975667346491, 0.00000005560486515676644, "Sí"),
    @("MCPS", 3, 7.135638097145171, 0.000000000005720979245893432, "Sí"),
    @("MCPS", 4, 7.904473336526208, 0.00000000000002375877272697835, "Sí"),
    @("MCPS", 5, 7.50142146785267, 0.0000000000004505285033928885, "Sí"),
    @("MCPS", 6, 7.268537000986087, 0.0000000000023046009545169, "Sí"),
    @("MCPS", 7, 9.249954030775257, 0, "Sí"),
    @("MCPS", 10, 7.874343128275965, 0.0000000000000297539770599542, "Sí"),
    @("Sieve Bootstrap", 1, 0.5333719674810627, 0.6067316732500676, "No"),
    @("Sieve Bootstrap", 2, 0.8679322243904682, 0.4022439599054577, "No"),
    @("Sieve Bootstrap", 3, -0.01864623655656174, 0.9856428175419636, "No"),
    @("Sieve Bootstrap", 4, 0.1215608505039952, 0.9066100980529868, "No"),
    @("Sieve Bootstrap", 5, 1.894091908761226, 0.06755822882963347, "No"),
    @("Sieve Bootstrap", 6, 6.94495170260025, 0.00000000002050182246193799, "Sí"),
    @("Sieve Bootstrap", 7, 6.577379263357245, 0.0000000002186013592364588, "Sí"),
    @("Sieve Bootstrap", 10, 7.34435645151481, 0.000000000001362021606610142, "Sí")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

Write-Host "Wrote $($data.Length) rows"
